$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the newest newsletter link (January 2019) in the row below the last entry
$ws.Range("B6").Value = "https://myemail.constantcontact.com/News-From-The-Forest--January-2019.html?soid=1102494320279&aid=GtYWC4C0xkA"

# Register it as a real hyperlink (matches the style used by the other links)
$ws.Hyperlinks.Add($ws.Range("B6"), "https://myemail.constantcontact.com/News-From-The-Forest--January-2019.html?soid=1102494320279&aid=GtYWC4C0xkA")

# Re-apply the same formatting (Hyperlink style) used by the cell above so B6 looks identical to B5
$ws.Range("B5").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move/keep the active selection where the user left off after adding the new row
$ws.Range("B17").Select()
